# demo-02-advanced.xlsx: add a third demo worksheet ("Demo3") that shows
# off the new "force fill rows" behaviour (setReadArea()/KEYS_FIRST_ROW) -
# a small staircase-shaped table with a header row styled like the
# existing Demo1/Demo2 header rows.

$wb = $excel.ActiveWorkbook

# Demo2 holds the header style (bold, centered, bordered) we want to reuse
# for Demo3's header row.
$ws2 = $wb.Worksheets.Item("Demo2")

# Add the new sheet right after the last existing sheet so the tab order
# ends up Demo1, Demo2, Demo3 - and it becomes the active sheet, same as
# in the target workbook (activeTab points at it, tabSelected moves to
# its sheetView).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Demo3"

# Header row (row 2): four new shared strings aaa/bbb/ccc/ddd.
$ws3.Range("A2").Value = "aaa"
$ws3.Range("B2").Value = "bbb"
$ws3.Range("C2").Value = "ccc"
$ws3.Range("D2").Value = "ddd"

# Copy the bold/centered/bordered header formatting from Demo2!B4 onto the
# new header row (this is the same style used for every other header row
# in the workbook).
$ws2.Range("B4").Copy()
$ws3.Range("A2:D2").PasteSpecial(-4122)

# Data rows 3-6: a "staircase" of numbers, one fewer value per row, to
# exercise the new force-fill-rows behaviour.
$ws3.Range("A3").Value = 3
$ws3.Range("B3").Value = 33
$ws3.Range("C3").Value = 333
$ws3.Range("D3").Value = 3333

$ws3.Range("A4").Value = 4
$ws3.Range("B4").Value = 44
$ws3.Range("C4").Value = 444

$ws3.Range("A5").Value = 5
$ws3.Range("B5").Value = 55

$ws3.Range("A6").Value = 6
